{"js": "/*\n * Office.js (Word JavaScript API) script implementing the commit:\n *   \"more design doc updates\" \u2014 adds a bold attribute to the\n *   \"Programming Approach\" heading, and appends two new subsections\n *   (\"Auxiliary Functions\" and \"ISRs (Interrupt Service Routines)\")\n *   right after the \"Illuminate corresponding LEDs.\" bullet, before\n *   the \"Wiring Schematic\" section.\n */\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"text\");\nawait context.sync();\n\n// --- 1) Bold the \"Programming Approach\" heading (pPr/rPr + run rPr both\n//        gain <w:b/>, matching the diff). ---------------------------------\nlet headingIndex = -1;\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  if (paragraphs.items[i].text === \"Programming Approach\") {\n    headingIndex = i;\n    break;\n  }\n}\nif (headingIndex === -1) {\n  throw new Error('Could not locate the \"Programming Approach\" heading paragraph.');\n}\nparagraphs.items[headingIndex].font.bold = true;\n\n// --- 2) Insert the new \"Auxiliary Functions\" / \"ISRs\" content -------------\n// Locate the \"\u2022 Illuminate corresponding LEDs.\" bullet paragraph; the new\n// content is inserted immediately after it (and before the pre-existing\n// blank paragraph that leads into the \"Wiring Schematic\" section).\nlet targetIndex = -1;\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  if (paragraphs.items[i].text === \"\\u2022 Illuminate corresponding LEDs.\") {\n    targetIndex = i;\n    break;\n  }\n}\nif (targetIndex === -1) {\n  throw new Error('Could not locate the \"Illuminate corresponding LEDs.\" paragraph.');\n}\n\nconst insertionRange = paragraphs.items[targetIndex].getRange(\"End\");\n\n// The new paragraphs, expressed as raw WordprocessingML, so formatting\n// details (fonts, bold/italic/size, proofErr spell-check ranges, and the\n// _GoBack bookmark) come through exactly as authored.\nconst newParagraphsXml = \"    <w:p>\\n      <w:pPr>\\n        <w:ind w:left=\\\"720\\\"/>\\n        <w:rPr>\\n          <w:rFonts w:ascii=\\\"Calibri\\\" w:hAnsi=\\\"Calibri\\\" w:cs=\\\"Menlo\\\"/>\\n        </w:rPr>\\n      </w:pPr>\\n    </w:p>\\n    <w:p>\\n      <w:pPr>\\n        <w:rPr>\\n          <w:rFonts w:ascii=\\\"Calibri\\\" w:hAnsi=\\\"Calibri\\\" w:cs=\\\"Menlo\\\"/>\\n          <w:b/>\\n          <w:i/>\\n          <w:sz w:val=\\\"28\\\"/>\\n        </w:rPr>\\n      </w:pPr>\\n      <w:r>\\n        <w:rPr>\\n          <w:rFonts w:ascii=\\\"Calibri\\\" w:hAnsi=\\\"Calibri\\\" w:cs=\\\"Menlo\\\"/>\\n          <w:b/>\\n          <w:i/>\\n          <w:sz w:val=\\\"28\\\"/>\\n        </w:rPr>\\n        <w:lastRenderedPageBreak/>\\n        <w:t>Auxiliary Functions</w:t>\\n      </w:r>\\n    </w:p>\\n    <w:p>\\n      <w:pPr>\\n        <w:spacing w:line=\\\"360\\\" w:lineRule=\\\"auto\\\"/>\\n        <w:rPr>\\n          <w:rFonts w:ascii=\\\"Menlo\\\" w:hAnsi=\\\"Menlo\\\" w:cs=\\\"Menlo\\\"/>\\n          <w:sz w:val=\\\"22\\\"/>\\n        </w:rPr>\\n      </w:pPr>\\n      <w:r>\\n        <w:rPr>\\n          <w:rFonts w:ascii=\\\"Menlo\\\" w:hAnsi=\\\"Menlo\\\" w:cs=\\\"Menlo\\\"/>\\n          <w:sz w:val=\\\"22\\\"/>\\n        </w:rPr>\\n        <w:t xml:space=\\\"preserve\\\">void </w:t>\\n      </w:r>\\n      <w:proofErr w:type=\\\"spellStart\\\"/>\\n      <w:r>\\n        <w:rPr>\\n          <w:rFonts w:ascii=\\\"Menlo\\\" w:hAnsi=\\\"Menlo\\\" w:cs=\\\"Menlo\\\"/>\\n          <w:sz w:val=\\\"22\\\"/>\\n        </w:rPr>\\n        <w:t>send_</w:t>\\n      </w:r>\\n      <w:proofErr w:type=\\\"gramStart\\\"/>\\n      <w:r>\\n        <w:rPr>\\n          <w:rFonts w:ascii=\\\"Menlo\\\" w:hAnsi=\\\"Menlo\\\" w:cs=\\\"Menlo\\\"/>\\n          <w:sz w:val=\\\"22\\\"/>\\n        </w:rPr>\\n        <w:t>timestamp</w:t>\\n      </w:r>\\n      <w:proofErr w:type=\\\"spellEnd\\\"/>\\n      <w:r>\\n        <w:rPr>\\n          <w:rFonts w:ascii=\\\"Menlo\\\" w:hAnsi=\\\"Menlo\\\" w:cs=\\\"Menlo\\\"/>\\n          <w:sz w:val=\\\"22\\\"/>\\n        </w:rPr>\\n        <w:t>(</w:t>\\n      </w:r>\\n      <w:proofErr w:type=\\\"gramEnd\\\"/>\\n      <w:r>\\n        <w:rPr>\\n          <w:rFonts w:ascii=\\\"Menlo\\\" w:hAnsi=\\\"Menlo\\\" w:cs=\\\"Menlo\\\"/>\\n          <w:sz w:val=\\\"22\\\"/>\\n        </w:rPr>\\n        <w:t>)</w:t>\\n      </w:r>\\n    </w:p>\\n    <w:p>\\n      <w:pPr>\\n        <w:spacing w:line=\\\"360\\\" w:lineRule=\\\"auto\\\"/>\\n        <w:rPr>\\n          <w:rFonts w:ascii=\\\"Calibri\\\" w:hAnsi=\\\"Calibri\\\" w:cs=\\\"Menlo\\\"/>\\n        </w:rPr>\\n      </w:pPr>\\n      <w:r>\\n        <w:rPr>\\n          <w:rFonts w:ascii=\\\"Calibri\\\" w:hAnsi=\\\"Calibri\\\" w:cs=\\\"Menlo\\\"/>\\n        </w:rPr>\\n        <w:t>This</w:t>\\n      </w:r>\\n      <w:r>\\n        <w:rPr>\\n          <w:rFonts w:ascii=\\\"Calibri\\\" w:hAnsi=\\\"Calibri\\\" w:cs=\\\"Menlo\\\"/>\\n        </w:rPr>\\n        <w:t xml:space=\\\"preserve\\\"> </w:t>\\n      </w:r>\\n      <w:r>\\n        <w:rPr>\\n          <w:rFonts w:ascii=\\\"Calibri\\\" w:hAnsi=\\\"Calibri\\\" w:cs=\\\"Menlo\\\"/>\\n        </w:rPr>\\n        <w:t xml:space=\\\"preserve\\\">is an auxiliary logger function that takes in 0 arguments, that logs the activity of the fan component of the system. It logs the fan\\u2019s status (on/off), and the </w:t>\\n      </w:r>\\n      <w:proofErr w:type=\\\"spellStart\\\"/>\\n      <w:r>\\n        <w:rPr>\\n          <w:rFonts w:ascii=\\\"Calibri\\\" w:hAnsi=\\\"Calibri\\\" w:cs=\\\"Menlo\\\"/>\\n        </w:rPr>\\n        <w:t>datetime</w:t>\\n      </w:r>\\n      <w:proofErr w:type=\\\"spellEnd\\\"/>\\n      <w:r>\\n        <w:rPr>\\n          <w:rFonts w:ascii=\\\"Calibri\\\" w:hAnsi=\\\"Calibri\\\" w:cs=\\\"Menlo\\\"/>\\n        </w:rPr>\\n        <w:t>. This is important not only because it is a basic system requirement of the project, but almost every single real-world application of embedded systems will have some sort of logging.</w:t>\\n      </w:r>\\n    </w:p>\\n    <w:p>\\n      <w:pPr>\\n        <w:rPr>\\n          <w:rFonts w:ascii=\\\"Calibri\\\" w:hAnsi=\\\"Calibri\\\" w:cs=\\\"Menlo\\\"/>\\n        </w:rPr>\\n      </w:pPr>\\n    </w:p>\\n    <w:p>\\n      <w:pPr>\\n        <w:rPr>\\n          <w:rFonts w:ascii=\\\"Calibri\\\" w:hAnsi=\\\"Calibri\\\" w:cs=\\\"Menlo\\\"/>\\n          <w:b/>\\n          <w:i/>\\n          <w:sz w:val=\\\"28\\\"/>\\n        </w:rPr>\\n      </w:pPr>\\n      <w:r>\\n        <w:rPr>\\n          <w:rFonts w:ascii=\\\"Calibri\\\" w:hAnsi=\\\"Calibri\\\" w:cs=\\\"Menlo\\\"/>\\n          <w:b/>\\n          <w:i/>\\n          <w:sz w:val=\\\"28\\\"/>\\n        </w:rPr>\\n        <w:t>ISRs (Interrupt Service Routines)</w:t>\\n      </w:r>\\n    </w:p>\\n    <w:p>\\n      <w:pPr>\\n        <w:spacing w:line=\\\"360\\\" w:lineRule=\\\"auto\\\"/>\\n        <w:rPr>\\n          <w:rFonts w:ascii=\\\"Menlo\\\" w:hAnsi=\\\"Menlo\\\" w:cs=\\\"Menlo\\\"/>\\n          <w:sz w:val=\\\"22\\\"/>\\n          <w:szCs w:val=\\\"22\\\"/>\\n        </w:rPr>\\n      </w:pPr>\\n      <w:r>\\n        <w:rPr>\\n          <w:rFonts w:ascii=\\\"Menlo\\\" w:hAnsi=\\\"Menlo\\\" w:cs=\\\"Menlo\\\"/>\\n          <w:sz w:val=\\\"22\\\"/>\\n          <w:szCs w:val=\\\"22\\\"/>\\n        </w:rPr>\\n        <w:t>ISR (TIMER1_OVF_vect)</w:t>\\n      </w:r>\\n    </w:p>\\n    <w:p>\\n      <w:pPr>\\n        <w:spacing w:line=\\\"360\\\" w:lineRule=\\\"auto\\\"/>\\n        <w:rPr>\\n          <w:rFonts w:ascii=\\\"Calibri\\\" w:hAnsi=\\\"Calibri\\\" w:cs=\\\"Menlo\\\"/>\\n          <w:szCs w:val=\\\"22\\\"/>\\n        </w:rPr>\\n      </w:pPr>\\n      <w:r>\\n        <w:rPr>\\n          <w:rFonts w:ascii=\\\"Calibri\\\" w:hAnsi=\\\"Calibri\\\" w:cs=\\\"Menlo\\\"/>\\n          <w:szCs w:val=\\\"22\\\"/>\\n        </w:rPr>\\n        <w:t xml:space=\\\"preserve\\\">This ISR\\u2019s purpose is to update the values of the temperature and humidity, with the frequency depending on the amount of ticks that is specified for </w:t>\\n      </w:r>\\n      <w:proofErr w:type=\\\"spellStart\\\"/>\\n      <w:r>\\n        <w:rPr>\\n          <w:rFonts w:ascii=\\\"Menlo\\\" w:hAnsi=\\\"Menlo\\\" w:cs=\\\"Menlo\\\"/>\\n          <w:sz w:val=\\\"22\\\"/>\\n          <w:szCs w:val=\\\"22\\\"/>\\n        </w:rPr>\\n        <w:t>temperature_humidity_sensor_sampling_tick</w:t>\\n      </w:r>\\n      <w:proofErr w:type=\\\"spellEnd\\\"/>\\n      <w:r>\\n        <w:rPr>\\n          <w:rFonts w:ascii=\\\"Calibri\\\" w:hAnsi=\\\"Calibri\\\" w:cs=\\\"Menlo\\\"/>\\n          <w:sz w:val=\\\"22\\\"/>\\n          <w:szCs w:val=\\\"22\\\"/>\\n        </w:rPr>\\n        <w:t xml:space=\\\"preserve\\\"> </w:t>\\n      </w:r>\\n      <w:r>\\n        <w:rPr>\\n          <w:rFonts w:ascii=\\\"Calibri\\\" w:hAnsi=\\\"Calibri\\\" w:cs=\\\"Menlo\\\"/>\\n          <w:szCs w:val=\\\"22\\\"/>\\n        </w:rPr>\\n        <w:t>macro.</w:t>\\n      </w:r>\\n    </w:p>\\n    <w:p>\\n      <w:pPr>\\n        <w:spacing w:line=\\\"360\\\" w:lineRule=\\\"auto\\\"/>\\n        <w:rPr>\\n          <w:rFonts w:ascii=\\\"Calibri\\\" w:hAnsi=\\\"Calibri\\\" w:cs=\\\"Menlo\\\"/>\\n          <w:szCs w:val=\\\"22\\\"/>\\n        </w:rPr>\\n      </w:pPr>\\n    </w:p>\\n    <w:p>\\n      <w:pPr>\\n        <w:spacing w:line=\\\"360\\\" w:lineRule=\\\"auto\\\"/>\\n        <w:rPr>\\n          <w:rFonts w:ascii=\\\"Menlo\\\" w:hAnsi=\\\"Menlo\\\" w:cs=\\\"Menlo\\\"/>\\n          <w:sz w:val=\\\"22\\\"/>\\n          <w:szCs w:val=\\\"22\\\"/>\\n        </w:rPr>\\n      </w:pPr>\\n      <w:r>\\n        <w:rPr>\\n          <w:rFonts w:ascii=\\\"Menlo\\\" w:hAnsi=\\\"Menlo\\\" w:cs=\\\"Menlo\\\"/>\\n          <w:sz w:val=\\\"22\\\"/>\\n          <w:szCs w:val=\\\"22\\\"/>\\n        </w:rPr>\\n        <w:t>ISR (INT3_vect)</w:t>\\n      </w:r>\\n    </w:p>\\n    <w:p>\\n      <w:pPr>\\n        <w:spacing w:line=\\\"360\\\" w:lineRule=\\\"auto\\\"/>\\n        <w:rPr>\\n          <w:rFonts w:ascii=\\\"Calibri\\\" w:hAnsi=\\\"Calibri\\\" w:cs=\\\"Menlo\\\"/>\\n          <w:szCs w:val=\\\"22\\\"/>\\n        </w:rPr>\\n      </w:pPr>\\n      <w:r>\\n        <w:rPr>\\n          <w:rFonts w:ascii=\\\"Calibri\\\" w:hAnsi=\\\"Calibri\\\" w:cs=\\\"Menlo\\\"/>\\n          <w:szCs w:val=\\\"22\\\"/>\\n        </w:rPr>\\n        <w:t xml:space=\\\"preserve\\\">This ISR\\u2019s purpose is to handle the </w:t>\\n      </w:r>\\n      <w:proofErr w:type=\\\"spellStart\\\"/>\\n      <w:r>\\n        <w:rPr>\\n          <w:rFonts w:ascii=\\\"Calibri\\\" w:hAnsi=\\\"Calibri\\\" w:cs=\\\"Menlo\\\"/>\\n          <w:szCs w:val=\\\"22\\\"/>\\n        </w:rPr>\\n        <w:t>debouncing</w:t>\\n      </w:r>\\n      <w:proofErr w:type=\\\"spellEnd\\\"/>\\n      <w:r>\\n        <w:rPr>\\n          <w:rFonts w:ascii=\\\"Calibri\\\" w:hAnsi=\\\"Calibri\\\" w:cs=\\\"Menlo\\\"/>\\n          <w:szCs w:val=\\\"22\\\"/>\\n        </w:rPr>\\n        <w:t xml:space=\\\"preserve\\\"> that occurs when/if one presses the pushbutton to force the system\\u2019s state to the disabled state.</w:t>\\n      </w:r>\\n    </w:p>\\n    <w:p>\\n      <w:pPr>\\n        <w:spacing w:line=\\\"360\\\" w:lineRule=\\\"auto\\\"/>\\n        <w:rPr>\\n          <w:rFonts w:ascii=\\\"Calibri\\\" w:hAnsi=\\\"Calibri\\\" w:cs=\\\"Menlo\\\"/>\\n          <w:szCs w:val=\\\"22\\\"/>\\n        </w:rPr>\\n      </w:pPr>\\n    </w:p>\\n    <w:p>\\n      <w:pPr>\\n        <w:spacing w:line=\\\"360\\\" w:lineRule=\\\"auto\\\"/>\\n        <w:rPr>\\n          <w:rFonts w:ascii=\\\"Menlo\\\" w:hAnsi=\\\"Menlo\\\" w:cs=\\\"Menlo\\\"/>\\n          <w:sz w:val=\\\"22\\\"/>\\n          <w:szCs w:val=\\\"22\\\"/>\\n        </w:rPr>\\n      </w:pPr>\\n      <w:r>\\n        <w:rPr>\\n          <w:rFonts w:ascii=\\\"Menlo\\\" w:hAnsi=\\\"Menlo\\\" w:cs=\\\"Menlo\\\"/>\\n          <w:sz w:val=\\\"22\\\"/>\\n          <w:szCs w:val=\\\"22\\\"/>\\n        </w:rPr>\\n        <w:t>ISR (</w:t>\\n      </w:r>\\n      <w:proofErr w:type=\\\"spellStart\\\"/>\\n      <w:r>\\n        <w:rPr>\\n          <w:rFonts w:ascii=\\\"Menlo\\\" w:hAnsi=\\\"Menlo\\\" w:cs=\\\"Menlo\\\"/>\\n          <w:sz w:val=\\\"22\\\"/>\\n          <w:szCs w:val=\\\"22\\\"/>\\n        </w:rPr>\\n        <w:t>ADC_vect</w:t>\\n      </w:r>\\n      <w:proofErr w:type=\\\"spellEnd\\\"/>\\n      <w:r>\\n        <w:rPr>\\n          <w:rFonts w:ascii=\\\"Menlo\\\" w:hAnsi=\\\"Menlo\\\" w:cs=\\\"Menlo\\\"/>\\n          <w:sz w:val=\\\"22\\\"/>\\n          <w:szCs w:val=\\\"22\\\"/>\\n        </w:rPr>\\n        <w:t>)</w:t>\\n      </w:r>\\n    </w:p>\\n    <w:p>\\n      <w:pPr>\\n        <w:spacing w:line=\\\"360\\\" w:lineRule=\\\"auto\\\"/>\\n        <w:rPr>\\n          <w:rFonts w:ascii=\\\"Calibri\\\" w:hAnsi=\\\"Calibri\\\" w:cs=\\\"Menlo\\\"/>\\n          <w:szCs w:val=\\\"22\\\"/>\\n        </w:rPr>\\n      </w:pPr>\\n      <w:r>\\n        <w:rPr>\\n          <w:rFonts w:ascii=\\\"Calibri\\\" w:hAnsi=\\\"Calibri\\\" w:cs=\\\"Menlo\\\"/>\\n          <w:szCs w:val=\\\"22\\\"/>\\n        </w:rPr>\\n        <w:t xml:space=\\\"preserve\\\">This ISR\\u2019s purpose is to </w:t>\\n      </w:r>\\n      <w:r>\\n        <w:rPr>\\n          <w:rFonts w:ascii=\\\"Calibri\\\" w:hAnsi=\\\"Calibri\\\" w:cs=\\\"Menlo\\\"/>\\n          <w:szCs w:val=\\\"22\\\"/>\\n        </w:rPr>\\n        <w:t xml:space=\\\"preserve\\\">provide the value that </w:t>\\n      </w:r>\\n      <w:proofErr w:type=\\\"spellStart\\\"/>\\n      <w:r>\\n        <w:rPr>\\n          <w:rFonts w:ascii=\\\"Menlo\\\" w:hAnsi=\\\"Menlo\\\" w:cs=\\\"Menlo\\\"/>\\n          <w:sz w:val=\\\"22\\\"/>\\n          <w:szCs w:val=\\\"22\\\"/>\\n        </w:rPr>\\n        <w:t>water_level</w:t>\\n      </w:r>\\n      <w:proofErr w:type=\\\"spellEnd\\\"/>\\n      <w:r>\\n        <w:rPr>\\n          <w:rFonts w:ascii=\\\"Calibri\\\" w:hAnsi=\\\"Calibri\\\" w:cs=\\\"Menlo\\\"/>\\n          <w:sz w:val=\\\"22\\\"/>\\n          <w:szCs w:val=\\\"22\\\"/>\\n        </w:rPr>\\n        <w:t xml:space=\\\"preserve\\\"> </w:t>\\n      </w:r>\\n      <w:r>\\n        <w:rPr>\\n          <w:rFonts w:ascii=\\\"Calibri\\\" w:hAnsi=\\\"Calibri\\\" w:cs=\\\"Menlo\\\"/>\\n          <w:szCs w:val=\\\"22\\\"/>\\n        </w:rPr>\\n        <w:t>should be set to. It initiates an interrupt to the Stack, once the ADC has completed its conversion from reading the current water level.</w:t>\\n      </w:r>\\n      <w:r>\\n        <w:rPr>\\n          <w:rFonts w:ascii=\\\"Calibri\\\" w:hAnsi=\\\"Calibri\\\" w:cs=\\\"Menlo\\\"/>\\n          <w:szCs w:val=\\\"22\\\"/>\\n        </w:rPr>\\n        <w:t xml:space=\\\"preserve\\\"> occurring </w:t>\\n      </w:r>\\n      <w:bookmarkStart w:id=\\\"0\\\" w:name=\\\"_GoBack\\\"/>\\n      <w:bookmarkEnd w:id=\\\"0\\\"/>\\n    </w:p>\\n    <w:p>\\n      <w:pPr>\\n        <w:ind w:left=\\\"720\\\"/>\\n        <w:rPr>\\n          <w:rFonts w:ascii=\\\"Calibri\\\" w:hAnsi=\\\"Calibri\\\" w:cs=\\\"Menlo\\\"/>\\n        </w:rPr>\\n      </w:pPr>\\n    </w:p>\\n\";\n\nconst flatOpcXml = `<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?>\n<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">\n<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">\n<pkg:xmlData>\n<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\"><w:body>${newParagraphsXml}<w:sectPr/></w:body></w:document>\n</pkg:xmlData>\n</pkg:part>\n</pkg:package>`;\n\ninsertionRange.insertOoxml(flatOpcXml, Word.InsertLocation.after);\n\nawait context.sync();\n", "ps1": "# PowerShell / Word COM-interop script implementing the commit:\n#   \"more design doc updates\" - adds bold to the \"Programming Approach\"\n#   heading, and appends two new subsections (\"Auxiliary Functions\" and\n#   \"ISRs (Interrupt Service Routines)\") right after the\n#   \"Illuminate corresponding LEDs.\" bullet, before the blank paragraph\n#   that leads into \"Wiring Schematic\".\n\n$d = $word.ActiveDocument\n\n# --- 1) Bold the \"Programming Approach\" heading ----------------------------\n$headingIndex = -1\nfor ($i = 1; $i -le $d.Paragraphs.Count; $i++) {\n    $t = $d.Paragraphs.Item($i).Range.Text.Trim()\n    if ($t -eq \"Programming Approach\") {\n        $headingIndex = $i\n        break\n    }\n}\nif ($headingIndex -eq -1) {\n    throw \"Could not locate the 'Programming Approach' heading paragraph.\"\n}\n$d.Paragraphs.Item($headingIndex).Range.Bold = 1\n\n# --- 2) Insert the new \"Auxiliary Functions\" / \"ISRs\" content --------------\n$bullet = [char]8226\n$needle = \"$bullet Illuminate corresponding LEDs.\"\n$targetIndex = -1\nfor ($i = 1; $i -le $d.Paragraphs.Count; $i++) {\n    $t = $d.Paragraphs.Item($i).Range.Text.Trim()\n    if ($t -eq $needle) {\n        $targetIndex = $i\n        break\n    }\n}\nif ($targetIndex -eq -1) {\n    throw \"Could not locate the 'Illuminate corresponding LEDs.' paragraph.\"\n}\n\n# Make room with a fresh empty paragraph right after the target, then\n# replace that empty paragraph's contents with the raw WordprocessingML\n# for the new section (so fonts, bold/italic/size, proofErr spell-check\n# ranges, and the _GoBack bookmark all come through exactly as authored).\n$targetRange = $d.Paragraphs.Item($targetIndex).Range\n$targetRange.InsertParagraphAfter()\n$newParaIndex = $targetIndex + 1\n$newParaRange = $d.Paragraphs.Item($newParaIndex).Range\n\n$innerXml = @'\n    <w:p>\n      <w:pPr>\n        <w:ind w:left=\"720\"/>\n        <w:rPr>\n          <w:rFonts w:ascii=\"Calibri\" w:hAnsi=\"Calibri\" w:cs=\"Menlo\"/>\n        </w:rPr>\n      </w:pPr>\n    </w:p>\n    <w:p>\n      <w:pPr>\n        <w:rPr>\n          <w:rFonts w:ascii=\"Calibri\" w:hAnsi=\"Calibri\" w:cs=\"Menlo\"/>\n          <w:b/>\n          <w:i/>\n          <w:sz w:val=\"28\"/>\n        </w:rPr>\n      </w:pPr>\n      <w:r>\n        <w:rPr>\n          <w:rFonts w:ascii=\"Calibri\" w:hAnsi=\"Calibri\" w:cs=\"Menlo\"/>\n          <w:b/>\n          <w:i/>\n          <w:sz w:val=\"28\"/>\n        </w:rPr>\n        <w:lastRenderedPageBreak/>\n        <w:t>Auxiliary Functions</w:t>\n      </w:r>\n    </w:p>\n    <w:p>\n      <w:pPr>\n        <w:spacing w:line=\"360\" w:lineRule=\"auto\"/>\n        <w:rPr>\n          <w:rFonts w:ascii=\"Menlo\" w:hAnsi=\"Menlo\" w:cs=\"Menlo\"/>\n          <w:sz w:val=\"22\"/>\n        </w:rPr>\n      </w:pPr>\n      <w:r>\n        <w:rPr>\n          <w:rFonts w:ascii=\"Menlo\" w:hAnsi=\"Menlo\" w:cs=\"Menlo\"/>\n          <w:sz w:val=\"22\"/>\n        </w:rPr>\n        <w:t xml:space=\"preserve\">void </w:t>\n      </w:r>\n      <w:proofErr w:type=\"spellStart\"/>\n      <w:r>\n        <w:rPr>\n          <w:rFonts w:ascii=\"Menlo\" w:hAnsi=\"Menlo\" w:cs=\"Menlo\"/>\n          <w:sz w:val=\"22\"/>\n        </w:rPr>\n        <w:t>send_</w:t>\n      </w:r>\n      <w:proofErr w:type=\"gramStart\"/>\n      <w:r>\n        <w:rPr>\n          <w:rFonts w:ascii=\"Menlo\" w:hAnsi=\"Menlo\" w:cs=\"Menlo\"/>\n          <w:sz w:val=\"22\"/>\n        </w:rPr>\n        <w:t>timestamp</w:t>\n      </w:r>\n      <w:proofErr w:type=\"spellEnd\"/>\n      <w:r>\n        <w:rPr>\n          <w:rFonts w:ascii=\"Menlo\" w:hAnsi=\"Menlo\" w:cs=\"Menlo\"/>\n          <w:sz w:val=\"22\"/>\n        </w:rPr>\n        <w:t>(</w:t>\n      </w:r>\n      <w:proofErr w:type=\"gramEnd\"/>\n      <w:r>\n        <w:rPr>\n          <w:rFonts w:ascii=\"Menlo\" w:hAnsi=\"Menlo\" w:cs=\"Menlo\"/>\n          <w:sz w:val=\"22\"/>\n        </w:rPr>\n        <w:t>)</w:t>\n      </w:r>\n    </w:p>\n    <w:p>\n      <w:pPr>\n        <w:spacing w:line=\"360\" w:lineRule=\"auto\"/>\n        <w:rPr>\n          <w:rFonts w:ascii=\"Calibri\" w:hAnsi=\"Calibri\" w:cs=\"Menlo\"/>\n        </w:rPr>\n      </w:pPr>\n      <w:r>\n        <w:rPr>\n          <w:rFonts w:ascii=\"Calibri\" w:hAnsi=\"Calibri\" w:cs=\"Menlo\"/>\n        </w:rPr>\n        <w:t>This</w:t>\n      </w:r>\n      <w:r>\n        <w:rPr>\n          <w:rFonts w:ascii=\"Calibri\" w:hAnsi=\"Calibri\" w:cs=\"Menlo\"/>\n        </w:rPr>\n        <w:t xml:space=\"preserve\"> </w:t>\n      </w:r>\n      <w:r>\n        <w:rPr>\n          <w:rFonts w:ascii=\"Calibri\" w:hAnsi=\"Calibri\" w:cs=\"Menlo\"/>\n        </w:rPr>\n        <w:t xml:space=\"preserve\">is an auxiliary logger function that takes in 0 arguments, that logs the activity of the fan component of the system. It logs the fan\u2019s status (on/off), and the </w:t>\n      </w:r>\n      <w:proofErr w:type=\"spellStart\"/>\n      <w:r>\n        <w:rPr>\n          <w:rFonts w:ascii=\"Calibri\" w:hAnsi=\"Calibri\" w:cs=\"Menlo\"/>\n        </w:rPr>\n        <w:t>datetime</w:t>\n      </w:r>\n      <w:proofErr w:type=\"spellEnd\"/>\n      <w:r>\n        <w:rPr>\n          <w:rFonts w:ascii=\"Calibri\" w:hAnsi=\"Calibri\" w:cs=\"Menlo\"/>\n        </w:rPr>\n        <w:t>. This is important not only because it is a basic system requirement of the project, but almost every single real-world application of embedded systems will have some sort of logging.</w:t>\n      </w:r>\n    </w:p>\n    <w:p>\n      <w:pPr>\n        <w:rPr>\n          <w:rFonts w:ascii=\"Calibri\" w:hAnsi=\"Calibri\" w:cs=\"Menlo\"/>\n        </w:rPr>\n      </w:pPr>\n    </w:p>\n    <w:p>\n      <w:pPr>\n        <w:rPr>\n          <w:rFonts w:ascii=\"Calibri\" w:hAnsi=\"Calibri\" w:cs=\"Menlo\"/>\n          <w:b/>\n          <w:i/>\n          <w:sz w:val=\"28\"/>\n        </w:rPr>\n      </w:pPr>\n      <w:r>\n        <w:rPr>\n          <w:rFonts w:ascii=\"Calibri\" w:hAnsi=\"Calibri\" w:cs=\"Menlo\"/>\n          <w:b/>\n          <w:i/>\n          <w:sz w:val=\"28\"/>\n        </w:rPr>\n        <w:t>ISRs (Interrupt Service Routines)</w:t>\n      </w:r>\n    </w:p>\n    <w:p>\n      <w:pPr>\n        <w:spacing w:line=\"360\" w:lineRule=\"auto\"/>\n        <w:rPr>\n          <w:rFonts w:ascii=\"Menlo\" w:hAnsi=\"Menlo\" w:cs=\"Menlo\"/>\n          <w:sz w:val=\"22\"/>\n          <w:szCs w:val=\"22\"/>\n        </w:rPr>\n      </w:pPr>\n      <w:r>\n        <w:rPr>\n          <w:rFonts w:ascii=\"Menlo\" w:hAnsi=\"Menlo\" w:cs=\"Menlo\"/>\n          <w:sz w:val=\"22\"/>\n          <w:szCs w:val=\"22\"/>\n        </w:rPr>\n        <w:t>ISR (TIMER1_OVF_vect)</w:t>\n      </w:r>\n    </w:p>\n    <w:p>\n      <w:pPr>\n        <w:spacing w:line=\"360\" w:lineRule=\"auto\"/>\n        <w:rPr>\n          <w:rFonts w:ascii=\"Calibri\" w:hAnsi=\"Calibri\" w:cs=\"Menlo\"/>\n          <w:szCs w:val=\"22\"/>\n        </w:rPr>\n      </w:pPr>\n      <w:r>\n        <w:rPr>\n          <w:rFonts w:ascii=\"Calibri\" w:hAnsi=\"Calibri\" w:cs=\"Menlo\"/>\n          <w:szCs w:val=\"22\"/>\n        </w:rPr>\n        <w:t xml:space=\"preserve\">This ISR\u2019s purpose is to update the values of the temperature and humidity, with the frequency depending on the amount of ticks that is specified for </w:t>\n      </w:r>\n      <w:proofErr w:type=\"spellStart\"/>\n      <w:r>\n        <w:rPr>\n          <w:rFonts w:ascii=\"Menlo\" w:hAnsi=\"Menlo\" w:cs=\"Menlo\"/>\n          <w:sz w:val=\"22\"/>\n          <w:szCs w:val=\"22\"/>\n        </w:rPr>\n        <w:t>temperature_humidity_sensor_sampling_tick</w:t>\n      </w:r>\n      <w:proofErr w:type=\"spellEnd\"/>\n      <w:r>\n        <w:rPr>\n          <w:rFonts w:ascii=\"Calibri\" w:hAnsi=\"Calibri\" w:cs=\"Menlo\"/>\n          <w:sz w:val=\"22\"/>\n          <w:szCs w:val=\"22\"/>\n        </w:rPr>\n        <w:t xml:space=\"preserve\"> </w:t>\n      </w:r>\n      <w:r>\n        <w:rPr>\n          <w:rFonts w:ascii=\"Calibri\" w:hAnsi=\"Calibri\" w:cs=\"Menlo\"/>\n          <w:szCs w:val=\"22\"/>\n        </w:rPr>\n        <w:t>macro.</w:t>\n      </w:r>\n    </w:p>\n    <w:p>\n      <w:pPr>\n        <w:spacing w:line=\"360\" w:lineRule=\"auto\"/>\n        <w:rPr>\n          <w:rFonts w:ascii=\"Calibri\" w:hAnsi=\"Calibri\" w:cs=\"Menlo\"/>\n          <w:szCs w:val=\"22\"/>\n        </w:rPr>\n      </w:pPr>\n    </w:p>\n    <w:p>\n      <w:pPr>\n        <w:spacing w:line=\"360\" w:lineRule=\"auto\"/>\n        <w:rPr>\n          <w:rFonts w:ascii=\"Menlo\" w:hAnsi=\"Menlo\" w:cs=\"Menlo\"/>\n          <w:sz w:val=\"22\"/>\n          <w:szCs w:val=\"22\"/>\n        </w:rPr>\n      </w:pPr>\n      <w:r>\n        <w:rPr>\n          <w:rFonts w:ascii=\"Menlo\" w:hAnsi=\"Menlo\" w:cs=\"Menlo\"/>\n          <w:sz w:val=\"22\"/>\n          <w:szCs w:val=\"22\"/>\n        </w:rPr>\n        <w:t>ISR (INT3_vect)</w:t>\n      </w:r>\n    </w:p>\n    <w:p>\n      <w:pPr>\n        <w:spacing w:line=\"360\" w:lineRule=\"auto\"/>\n        <w:rPr>\n          <w:rFonts w:ascii=\"Calibri\" w:hAnsi=\"Calibri\" w:cs=\"Menlo\"/>\n          <w:szCs w:val=\"22\"/>\n        </w:rPr>\n      </w:pPr>\n      <w:r>\n        <w:rPr>\n          <w:rFonts w:ascii=\"Calibri\" w:hAnsi=\"Calibri\" w:cs=\"Menlo\"/>\n          <w:szCs w:val=\"22\"/>\n        </w:rPr>\n        <w:t xml:space=\"preserve\">This ISR\u2019s purpose is to handle the </w:t>\n      </w:r>\n      <w:proofErr w:type=\"spellStart\"/>\n      <w:r>\n        <w:rPr>\n          <w:rFonts w:ascii=\"Calibri\" w:hAnsi=\"Calibri\" w:cs=\"Menlo\"/>\n          <w:szCs w:val=\"22\"/>\n        </w:rPr>\n        <w:t>debouncing</w:t>\n      </w:r>\n      <w:proofErr w:type=\"spellEnd\"/>\n      <w:r>\n        <w:rPr>\n          <w:rFonts w:ascii=\"Calibri\" w:hAnsi=\"Calibri\" w:cs=\"Menlo\"/>\n          <w:szCs w:val=\"22\"/>\n        </w:rPr>\n        <w:t xml:space=\"preserve\"> that occurs when/if one presses the pushbutton to force the system\u2019s state to the disabled state.</w:t>\n      </w:r>\n    </w:p>\n    <w:p>\n      <w:pPr>\n        <w:spacing w:line=\"360\" w:lineRule=\"auto\"/>\n        <w:rPr>\n          <w:rFonts w:ascii=\"Calibri\" w:hAnsi=\"Calibri\" w:cs=\"Menlo\"/>\n          <w:szCs w:val=\"22\"/>\n        </w:rPr>\n      </w:pPr>\n    </w:p>\n    <w:p>\n      <w:pPr>\n        <w:spacing w:line=\"360\" w:lineRule=\"auto\"/>\n        <w:rPr>\n          <w:rFonts w:ascii=\"Menlo\" w:hAnsi=\"Menlo\" w:cs=\"Menlo\"/>\n          <w:sz w:val=\"22\"/>\n          <w:szCs w:val=\"22\"/>\n        </w:rPr>\n      </w:pPr>\n      <w:r>\n        <w:rPr>\n          <w:rFonts w:ascii=\"Menlo\" w:hAnsi=\"Menlo\" w:cs=\"Menlo\"/>\n          <w:sz w:val=\"22\"/>\n          <w:szCs w:val=\"22\"/>\n        </w:rPr>\n        <w:t>ISR (</w:t>\n      </w:r>\n      <w:proofErr w:type=\"spellStart\"/>\n      <w:r>\n        <w:rPr>\n          <w:rFonts w:ascii=\"Menlo\" w:hAnsi=\"Menlo\" w:cs=\"Menlo\"/>\n          <w:sz w:val=\"22\"/>\n          <w:szCs w:val=\"22\"/>\n        </w:rPr>\n        <w:t>ADC_vect</w:t>\n      </w:r>\n      <w:proofErr w:type=\"spellEnd\"/>\n      <w:r>\n        <w:rPr>\n          <w:rFonts w:ascii=\"Menlo\" w:hAnsi=\"Menlo\" w:cs=\"Menlo\"/>\n          <w:sz w:val=\"22\"/>\n          <w:szCs w:val=\"22\"/>\n        </w:rPr>\n        <w:t>)</w:t>\n      </w:r>\n    </w:p>\n    <w:p>\n      <w:pPr>\n        <w:spacing w:line=\"360\" w:lineRule=\"auto\"/>\n        <w:rPr>\n          <w:rFonts w:ascii=\"Calibri\" w:hAnsi=\"Calibri\" w:cs=\"Menlo\"/>\n          <w:szCs w:val=\"22\"/>\n        </w:rPr>\n      </w:pPr>\n      <w:r>\n        <w:rPr>\n          <w:rFonts w:ascii=\"Calibri\" w:hAnsi=\"Calibri\" w:cs=\"Menlo\"/>\n          <w:szCs w:val=\"22\"/>\n        </w:rPr>\n        <w:t xml:space=\"preserve\">This ISR\u2019s purpose is to </w:t>\n      </w:r>\n      <w:r>\n        <w:rPr>\n          <w:rFonts w:ascii=\"Calibri\" w:hAnsi=\"Calibri\" w:cs=\"Menlo\"/>\n          <w:szCs w:val=\"22\"/>\n        </w:rPr>\n        <w:t xml:space=\"preserve\">provide the value that </w:t>\n      </w:r>\n      <w:proofErr w:type=\"spellStart\"/>\n      <w:r>\n        <w:rPr>\n          <w:rFonts w:ascii=\"Menlo\" w:hAnsi=\"Menlo\" w:cs=\"Menlo\"/>\n          <w:sz w:val=\"22\"/>\n          <w:szCs w:val=\"22\"/>\n        </w:rPr>\n        <w:t>water_level</w:t>\n      </w:r>\n      <w:proofErr w:type=\"spellEnd\"/>\n      <w:r>\n        <w:rPr>\n          <w:rFonts w:ascii=\"Calibri\" w:hAnsi=\"Calibri\" w:cs=\"Menlo\"/>\n          <w:sz w:val=\"22\"/>\n          <w:szCs w:val=\"22\"/>\n        </w:rPr>\n        <w:t xml:space=\"preserve\"> </w:t>\n      </w:r>\n      <w:r>\n        <w:rPr>\n          <w:rFonts w:ascii=\"Calibri\" w:hAnsi=\"Calibri\" w:cs=\"Menlo\"/>\n          <w:szCs w:val=\"22\"/>\n        </w:rPr>\n        <w:t>should be set to. It initiates an interrupt to the Stack, once the ADC has completed its conversion from reading the current water level.</w:t>\n      </w:r>\n      <w:r>\n        <w:rPr>\n          <w:rFonts w:ascii=\"Calibri\" w:hAnsi=\"Calibri\" w:cs=\"Menlo\"/>\n          <w:szCs w:val=\"22\"/>\n        </w:rPr>\n        <w:t xml:space=\"preserve\"> occurring </w:t>\n      </w:r>\n      <w:bookmarkStart w:id=\"0\" w:name=\"_GoBack\"/>\n      <w:bookmarkEnd w:id=\"0\"/>\n    </w:p>\n    <w:p>\n      <w:pPr>\n        <w:ind w:left=\"720\"/>\n        <w:rPr>\n          <w:rFonts w:ascii=\"Calibri\" w:hAnsi=\"Calibri\" w:cs=\"Menlo\"/>\n        </w:rPr>\n      </w:pPr>\n    </w:p>\n\n'@\n\n$flatOpcXml = '<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?><pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\"><pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\"><pkg:xmlData><w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\"><w:body>' + $innerXml + '<w:sectPr/></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'\n\n$newParaRange.InsertXML($flatOpcXml)\n"}
